$d = $word.ActiveDocument
$t1 = $d.Tables.Item(1)

# Row: "Age at transplant, y" -> update Overall / Post Transplant / Pre Transplant columns
$t1.Cell(6, 2).Range.Text = "54 (12)"
$t1.Cell(6, 3).Range.Text = "49 (11)"
$t1.Cell(6, 4).Range.Text = "58 (11)"

# Row: "Age at RCC, y" -> update Overall / Post Transplant / Pre Transplant columns
$t1.Cell(7, 2).Range.Text = "56 (10)"
$t1.Cell(7, 3).Range.Text = "57 (10)"
$t1.Cell(7, 4).Range.Text = "55 (11)"

# Footnote legend row: "n (%); Median (Q1, Q3)" -> "n (%); Mean (SD)"
# (scope Find to this cell's range so the superscript "1" run is left untouched)
$t1.Cell(22, 2).Range.Find.Execute("n (%); Median (Q1, Q3)", $true, $false, $false, $false, $false, $true, 1, $false, "n (%); Mean (SD)", 2) | Out-Null
